$d = $word.ActiveDocument

# 1. Update the existing text of the last paragraph to add the parenthetical clause
$d.Content.Find.Execute(
    "something bad may happen.  The goal is clear",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "something bad may happen (or another way to say this is that the bird cannot be left with the other items).  The goal is clear",
    2
)

# 2. Append new paragraphs at the end of the document

# Empty paragraph
$r1 = $d.Content
$r1.Collapse(0)
$r1.InsertParagraphAfter()

# "Identify Potential solutions:"
$r2 = $d.Content
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$r2.Text = "Identify Potential solutions:"

# Long solutions paragraph
$r3 = $d.Content
$r3.Collapse(0)
$r3.InsertParagraphAfter()
$r3.Text = "First solution: The man will take the bird with him first.  He will then take the cat to the other side.  He will take the bird back with him to the starting side and leave it there as he takes the seed.  He will then go back finally for the bird.  Second solution (similar to the first but the order of the cat and seed transfer reversed): Bird first. Then take the seed (instead of the cat this time). Bring the bird back on the trip for the cat (instead of the seed). Leave the bird again and take the cat to the destination and leave it with the seed.  Finally go back for the bird."

# Trailing empty paragraph
$r4 = $d.Content
$r4.Collapse(0)
$r4.InsertParagraphAfter()

# Final trailing empty paragraph
$r5 = $d.Content
$r5.Collapse(0)
$r5.InsertParagraphAfter()
